# Update system inputs from Excel
# Insert three new battery-related input rows into the "gStation" sheet
# (batt.E_rated, batt.E_ex, batt.f_repl), shifting the existing
# hydAccum/hydMotor/pumpMotor rows down, then make "gStation" the active
# sheet/tab as it was left after the edit.

$wb = $excel.ActiveWorkbook

$gStation = $wb.Worksheets.Item("gStation")

# Make room for the new battery rows (old rows 4-8 shift to 7-11)
$gStation.Rows("4:6").Insert()

$gStation.Range("A4").Value = "batt.E_rated"
$gStation.Range("B4").Value = 1000

$gStation.Range("A5").Value = "batt.E_ex"
$gStation.Range("B5").Value = "[0	0	0	0	0	0.794997961209469	1.91620710592149	3.66340782590143	7.00053194254757	10.6616297528835	11.1890647627290	11.2279405187141	11.2277213125461	11.1934366085895	11.2473922571621	11.2098142616423	11.1438062076693	11.1273616907102	11.1228424486107	11.1240295153377	11.1278176039227	11.1328371206444	11.1385755705485	11.1448282168896	11.1514961274243]"

$gStation.Range("A6").Value = "batt.f_repl"
$gStation.Range("B6").Value = -1

# The workbook was left with "gStation" selected/active instead of "system"
$gStation.Activate()
$gStation.Range("B16").Select()
